# edit.ps1
# Applies the changes described by the commit "add fxn to assign manips to macd":
#  1. Fix typo "conduced" -> "conducted"
#  2. Clean up the diazrenata/feasiblesads run (merge into a single plain run)
#  3. Split the "For each focal SAD..." paragraph: keep only the first sentence
#  4. Insert a new paragraph discussing small feasible sets / exclusion criterion
#  5. Insert a new paragraph (moved content) about constructing skew/evenness
#     distributions and comparing SAD position
#  6. Insert a new sub-bullet "For now, try just the R2 of a lm()"

$d = $word.ActiveDocument
$apos = [char]8217   # Unicode right single quotation mark used in "SAD's"

# ---------------------------------------------------------------------------
# 1. Fix "conduced" -> "conducted"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("we conduced all analyses", $true, $false, $false, $false, $false, $true, 1, $false, "we conducted all analyses", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Clean up "diazrenata/feasiblesads" -- replace the whole run sequence
#    (which includes proofErr spell-check wrappers) with a single plain run
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute(" implemented in diazrenata/feasiblesads", $true, $false, $false, $false, $false, $true, 1, $false, " implemented in diazrenata/feasiblesads", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Trim the "For each focal SAD..." paragraph down to its first sentence,
#    removing the second sentence and the "This allows us..." sentences
#    (those will be reinstated, in part, in a new paragraph below).
# ---------------------------------------------------------------------------
$oldFocalText = "For each focal SAD, we drew 2500 samples from the feasible set of an abundance distribution with the appropriate S and N. We constructed the distribution of skewness and evenness values of the sampled SADs, and calculated the percentile rank of the focal SAD" + $apos + "s statistics relative to these distributions. This allows us to compare the position of SADs from communities with different S and N and therefore different feasible sets. "
$newFocalText = "For each focal SAD, we drew 2500 samples from the feasible set of an abundance distribution with the appropriate S and N. "
$rng = $d.Content
$rng.Find.Execute($oldFocalText, $true, $false, $false, $false, $false, $true, 1, $false, $newFocalText, 2) | Out-Null

# Locate this (now-trimmed) paragraph so we can insert new paragraphs after it.
$rng = $d.Content
$rng.Find.Execute($newFocalText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Expand(4) | Out-Null   # wdParagraph -> grab the whole paragraph

# ---------------------------------------------------------------------------
# 4. Insert new paragraph: "Some combinations of S and N ... subsequent analyses."
#    (same indent level / list style as the paragraph it follows)
# ---------------------------------------------------------------------------
$rng.InsertParagraphAfter() | Out-Null
$insPoint1 = $d.Range($rng.End, $rng.End)
$insPoint1.InsertAfter("Some combinations of S and N have a feasible set that is much smaller than 2500. It may not be appropriate to draw conclusions based on such a small feasible space. If 2500 draws from the feasible set did not yield 2000 or more unique distributions, we did not include that community in subsequent analyses. ") | Out-Null

# ---------------------------------------------------------------------------
# 5. Insert new paragraph with the previously-removed sentences (skew/evenness
#    distribution + comparing SAD position), with "position" italicised.
# ---------------------------------------------------------------------------
$rng2 = $d.Range($insPoint1.Start, $insPoint1.End)
$rng2.Expand(4) | Out-Null   # whole "Some combinations..." paragraph
$rng2.InsertParagraphAfter() | Out-Null

$sentence1 = "We constructed the distribution of skewness and evenness values of the sampled SADs, and calculated the percentile rank of the focal SAD" + $apos + "s statistics relative to these distributions. "
$phraseA = "This allows us to compare the "
$phraseItalic = "position "
$phraseB = "of SADs "
$phraseC = "from communities"
$phraseD = " with different S and N and therefore different feasible sets. "

$newPara2Start = $rng2.End
$insPoint2 = $d.Range($newPara2Start, $newPara2Start)
$insPoint2.InsertAfter($sentence1 + $phraseA + $phraseItalic + $phraseB + $phraseC + $phraseD) | Out-Null

# Italicise just the word "position " within the paragraph we just filled in.
$italicStart = $newPara2Start + ($sentence1 + $phraseA).Length
$italicEnd = $italicStart + $phraseItalic.Length
$italicRng = $d.Range($italicStart, $italicEnd)
$italicRng.Font.Italic = 1

# ---------------------------------------------------------------------------
# 6. Insert new sub-bullet "For now, try just the R2 of a lm()" after the
#    "Beta regression..." bullet.
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Beta regression to test whether average abundance, mean/sd of FS statistic predict percentile position?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng3.Expand(4) | Out-Null
$rng3.InsertParagraphAfter() | Out-Null
$insPoint3 = $d.Range($rng3.End, $rng3.End)
$insPoint3.InsertAfter("For now, try just the R2 of a lm()") | Out-Null

Write-Host "edit complete"
